$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.187.81"
$ws.Range("E2").Value = "  +1.56%  "

$ws.Range("D3").Value = "1.437.81"
$ws.Range("E3").Value = "  +3.26%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.65%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9146"
$ws.Range("E5").Value = "  -8.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "276.23"
$ws.Range("E6").Value = "  +2.63%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3623"
$ws.Range("E7").Value = "  -0.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3073"
$ws.Range("E8").Value = "  +0.38%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "38.88"
$ws.Range("E9").Value = "  -0.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.023"
$ws.Range("E10").Value = "  +3.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06485"
$ws.Range("E11").Value = "  +1.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9996"
$ws.Range("E12").Value = "  +0.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.339"
$ws.Range("E13").Value = "  +0.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.49"
$ws.Range("E14").Value = "  +3.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.035"
$ws.Range("E15").Value = "  -0.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001009"
$ws.Range("E16").Value = "  +1.17%  "

$ws.Range("D17").Value = "1.432.25"
$ws.Range("E17").Value = "  +2.94%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9396"
$ws.Range("E18").Value = "  -5.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05632"
$ws.Range("E19").Value = "  -0.33%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.60"
$ws.Range("E20").Value = "  -3.96%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.378"
$ws.Range("E21").Value = "  -2.56%  "

$ws.Range("E22").Value = "  -2.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.86"
$ws.Range("E23").Value = "  +1.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.243"
$ws.Range("E24").Value = "  +0.16%  "

$ws.Range("D25").Value = "20.184.42"
$ws.Range("E25").Value = "  +1.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.05"
$ws.Range("E26").Value = "  +2.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.127"
$ws.Range("E27").Value = "  -2.90%  "

$ws.Range("E28").Value = "  +1.20%  "

$ws.Range("D29").Value = "1.585.29"
$ws.Range("E29").Value = "  +2.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "110.00"
$ws.Range("E30").Value = "  +1.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.820"
$ws.Range("E31").Value = "  -6.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8083"
$ws.Range("E32").Value = "  +0.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.820"
$ws.Range("E33").Value = "  -7.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07637"
$ws.Range("E34").Value = "  +0.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.466"
$ws.Range("E35").Value = "  +9.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05820"
$ws.Range("E36").Value = "  +3.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.129"
$ws.Range("E37").Value = "  +6.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.618"
$ws.Range("E38").Value = "  -2.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01985"
$ws.Range("E39").Value = "  -2.41%  "

$ws.Range("E40").Value = "  -0.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9232"
$ws.Range("E41").Value = "  -7.48%  "

$ws.Range("E42").Value = "  -2.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.135"
$ws.Range("E43").Value = "  -13.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.494"
$ws.Range("E44").Value = "  +0.64%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5197"
$ws.Range("E45").Value = "  -0.33%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.78"
$ws.Range("E46").Value = "  -2.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "116.51"
$ws.Range("E47").Value = "  +5.94%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5071"
$ws.Range("E48").Value = "  +1.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.742"
$ws.Range("E49").Value = "  -0.17%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06334"
$ws.Range("E50").Value = "  +3.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9876"
$ws.Range("E51").Value = "  -1.05%  "
